$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" placeholder on the Slide Master and on every
#    Slide Layout: the cached date text changes from 8/25/2021 to 8/27/2021.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                if ($sh.TextFrame.TextRange.Text -eq "8/25/2021") {
                    $sh.TextFrame.TextRange.Text = "8/27/2021"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 13 ("How to locate a feature on a slope?"): the "Arrow: Right 9"
#    shape is nudged from (1358283, 4358936) EMU to (1362919, 4367813) EMU.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
for ($i = 1; $i -le $s13.Shapes.Count; $i++) {
    $sh = $s13.Shapes.Item($i)
    if ($sh.Name -eq "Right Arrow 9") {
        $sh.Left = 107.31645669291339
        $sh.Top = 343.92230224609375
    }
}
